$wb = $excel.ActiveWorkbook

# This edit refreshes market-price-derived columns (H..N) across several
# sheets, as produced by the scheduled profit-tracker runner. Values come
# from an external pricing source; there are no formulas to recompute.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1100
$ws.Range("J29").Value = 2500
$ws.Range("L29").Value = 7500
$ws.Range("N29").Value = -8062

$ws.Range("H31").Value = 72
$ws.Range("I31").Value = 72
$ws.Range("K31").Value = 216
$ws.Range("M31").Value = 14

$ws.Range("H38").Value = 952.2857
$ws.Range("I38").Value = 111
$ws.Range("J38").Value = 6000
$ws.Range("K38").Value = 333
$ws.Range("L38").Value = 18000
$ws.Range("M38").Value = 39
$ws.Range("N38").Value = -18744

$ws.Range("H87").Value = 57800
$ws.Range("J87").Value = 57800
$ws.Range("L87").Value = 57800
$ws.Range("N87").Value = -60296

$ws.Range("H90").Value = 57800
$ws.Range("J90").Value = 57800
$ws.Range("L90").Value = 173400
$ws.Range("N90").Value = -185880

$ws.Range("H138").Value = 2666.6667
$ws.Range("J138").Value = 2296.3333
$ws.Range("L138").Value = 6888.999899999999
$ws.Range("N138").Value = -17168.9999

$ws.Range("H141").Value = 2529.2942
$ws.Range("J141").Value = 2145
$ws.Range("L141").Value = 6435
$ws.Range("N141").Value = -16795

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 18658.334
$ws.Range("I37").Value = 17987.5
$ws.Range("J37").Value = 20000
$ws.Range("K37").Value = 17987.5
$ws.Range("L37").Value = 20000
$ws.Range("M37").Value = -17714.5
$ws.Range("N37").Value = -20546

$ws.Range("H55").Value = 21479.666
$ws.Range("J55").Value = 34489
$ws.Range("L55").Value = 34489
$ws.Range("N55").Value = -35119

$ws.Range("H61").Value = 5091.156
$ws.Range("I61").Value = 3944.225
$ws.Range("J61").Value = 14266.6
$ws.Range("K61").Value = 3944.225
$ws.Range("L61").Value = 14266.6
$ws.Range("M61").Value = -3732.225
$ws.Range("N61").Value = -14690.6

$ws.Range("H102").Value = 4117.0415
$ws.Range("I102").Value = 2990.9524
$ws.Range("K102").Value = 2990.9524
$ws.Range("M102").Value = -1368.9524

$ws.Range("H136").Value = 5091.156
$ws.Range("I136").Value = 3944.225
$ws.Range("J136").Value = 14266.6
$ws.Range("K136").Value = 11832.675
$ws.Range("L136").Value = 42799.8
$ws.Range("M136").Value = -9282.674999999999
$ws.Range("N136").Value = -47899.8

$ws.Range("H137").Value = 89997.664
$ws.Range("J137").Value = 89997.664
$ws.Range("L137").Value = 89997.664
$ws.Range("N137").Value = -100197.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 31238
$ws.Range("J35").Value = 31238
$ws.Range("L35").Value = 31238
$ws.Range("N35").Value = -31858

$ws.Range("H82").Value = 23350
$ws.Range("J82").Value = 44500
$ws.Range("L82").Value = 44500
$ws.Range("N82").Value = -45266

$ws.Range("H85").Value = 23350
$ws.Range("J85").Value = 44500
$ws.Range("L85").Value = 44500
$ws.Range("N85").Value = -47152

$ws.Range("H94").Value = 5150.25
$ws.Range("I94").Value = 4600.2856
$ws.Range("J94").Value = 9000
$ws.Range("K94").Value = 4600.2856
$ws.Range("L94").Value = 9000
$ws.Range("M94").Value = -4149.2856
$ws.Range("N94").Value = -9902

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 3725.8
$ws.Range("I35").Value = 2658.5
$ws.Range("J35").Value = 7995
$ws.Range("K35").Value = 2658.5
$ws.Range("L35").Value = 7995
$ws.Range("M35").Value = -2364.5
$ws.Range("N35").Value = -8583

$ws.Range("H41").Value = 21923.572
$ws.Range("J41").Value = 21013.4
$ws.Range("L41").Value = 21013.4
$ws.Range("N41").Value = -21869.4

$ws.Range("H51").Value = 27800
$ws.Range("J51").Value = 27800
$ws.Range("L51").Value = 27800
$ws.Range("N51").Value = -29272

$ws.Range("H58").Value = 9272.182000000001
$ws.Range("J58").Value = 14999.8
$ws.Range("L58").Value = 14999.8
$ws.Range("N58").Value = -15405.8

$ws.Range("H59").Value = 40979
$ws.Range("J59").Value = 40979
$ws.Range("L59").Value = 40979
$ws.Range("N59").Value = -43269

$ws.Range("H60").Value = 27400
$ws.Range("J60").Value = 27400
$ws.Range("L60").Value = 27400
$ws.Range("N60").Value = -28422

$ws.Range("H61").Value = 27800
$ws.Range("J61").Value = 27800
$ws.Range("L61").Value = 27800
$ws.Range("N61").Value = -28496

$ws.Range("H68").Value = 41254.2
$ws.Range("J68").Value = 41254.2
$ws.Range("L68").Value = 41254.2
$ws.Range("N68").Value = -42752.2

$ws.Range("H71").Value = 41254.2
$ws.Range("J71").Value = 41254.2
$ws.Range("L71").Value = 123762.6
$ws.Range("N71").Value = -131250.6

$ws.Range("H74").Value = 37113.145
$ws.Range("J74").Value = 37113.145
$ws.Range("L74").Value = 37113.145
$ws.Range("N74").Value = -38861.145

$ws.Range("H77").Value = 37113.145
$ws.Range("J77").Value = 37113.145
$ws.Range("L77").Value = 111339.435
$ws.Range("N77").Value = -120075.435

$ws.Range("H134").Value = 4122.4546
$ws.Range("I134").Value = 3138.5454
$ws.Range("K134").Value = 9415.636200000001
$ws.Range("M134").Value = -6880.636200000001

$ws.Range("H136").Value = 9272.182000000001
$ws.Range("J136").Value = 14999.8
$ws.Range("L136").Value = 44999.39999999999
$ws.Range("N136").Value = -50099.39999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 8573.25
$ws.Range("I2").Value = 213
$ws.Range("K2").Value = 213
$ws.Range("M2").Value = -100

$ws.Range("H46").Value = 22953.092
$ws.Range("I46").Value = 15624.75
$ws.Range("J46").Value = 27140.715
$ws.Range("K46").Value = 15624.75
$ws.Range("L46").Value = 27140.715
$ws.Range("M46").Value = -15468.75
$ws.Range("N46").Value = -27452.715

$ws.Range("H51").Value = 41428.5
$ws.Range("I51").Value = 41428.5
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 41428.5
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -40919.5
$ws.Range("N51").ClearContents()

$ws.Range("H113").Value = 2782.8333
$ws.Range("I113").Value = 3474.25
$ws.Range("K113").Value = 3474.25
$ws.Range("M113").Value = -1304.25

$ws.Range("H132").Value = 2550.1738
$ws.Range("I132").Value = 2177.7104
$ws.Range("K132").Value = 6533.1312
$ws.Range("M132").Value = -4003.1312

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2012.7858
$ws.Range("I16").Value = 1709.3334
$ws.Range("K16").Value = 1709.3334
$ws.Range("M16").Value = -1539.3334

$ws.Range("H22").Value = 3911.182
$ws.Range("I22").Value = 3007
$ws.Range("K22").Value = 3007
$ws.Range("M22").Value = -2712

$ws.Range("H27").Value = 3911.182
$ws.Range("I27").Value = 3007
$ws.Range("K27").Value = 3007
$ws.Range("M27").Value = -2900

$ws.Range("H32").Value = 8318.5
$ws.Range("I32").Value = 8318.5
$ws.Range("K32").Value = 8318.5
$ws.Range("M32").Value = -8001.5

$ws.Range("H100").Value = 5275.3335
$ws.Range("I100").Value = 5826
$ws.Range("K100").Value = 5826
$ws.Range("M100").Value = -5285

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 5000250
$ws.Range("I3").Value = 5000250
$ws.Range("K3").Value = 5000250
$ws.Range("M3").Value = -5000136

$ws.Range("H41").Value = 10956.692
$ws.Range("I41").Value = 9000
$ws.Range("J41").Value = 11119.75
$ws.Range("K41").Value = 9000
$ws.Range("L41").Value = 11119.75
$ws.Range("M41").Value = -8610
$ws.Range("N41").Value = -11899.75

$ws.Range("H54").Value = 27600
$ws.Range("J54").Value = 27600
$ws.Range("L54").Value = 27600
$ws.Range("N54").Value = -28640

$ws.Range("H100").Value = 1184.4375
$ws.Range("I100").Value = 1165.9231
$ws.Range("K100").Value = 2331.8462
$ws.Range("M100").Value = -1790.8462

$ws.Range("H126").Value = 3159.6
$ws.Range("I126").Value = 3199.5
$ws.Range("K126").Value = 9598.5
$ws.Range("M126").Value = -7128.5

$ws.Range("H132").Value = 2322.8684
$ws.Range("I132").Value = 2224.4856
$ws.Range("K132").Value = 6673.4568
$ws.Range("M132").Value = -4143.4568
